# Meetings.xlsx - "Add files via upload" edit
# Adds per-meeting attendance (F:I) and duration (K) columns, plus two new
# discussion notes (C23, C25), for the meeting rows at r=23,25,27,29,31,33.
# Shared strings for the new duration values + discussion notes are created
# in the same order Excel would assign them (K23,K25,K27,K29,K33 durations,
# then C23,C25 notes) so the sharedStrings table layout matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New duration values in column K (also attendance columns F:I reuse
#     the existing "1)abhinav / 2)Harsha / 3)Udhay / 4)Prakyath" strings) ---
$ws.Range("K23").Value = "10:15-11:22"
$ws.Range("K25").Value = "10:15-11:24"
$ws.Range("K27").Value = "10:15-11:26"
$ws.Range("K29").Value = "10:15-11:28"
$ws.Range("K33").Value = "10:15-11:32"

# --- New discussion notes in column C ---
$ws.Range("C23").Value = "Demo was given on this date"
$ws.Range("C25").Value = "Code was inspected on this day"

# --- Attendance columns F:I for every meeting row, and the remaining
#     K31 duration which reuses the pre-existing "10:15-11:30" string ---
$rows = 23, 25, 27, 29, 31, 33
foreach ($r in $rows) {
    $ws.Range("F$r").Value = "1)abhinav"
    $ws.Range("G$r").Value = "2)Harsha"
    $ws.Range("H$r").Value = "3)Udhay"
    $ws.Range("I$r").Value = "4)Prakyath"
}
$ws.Range("K31").Value = "10:15-11:30"

# --- View state: last selected cell ---
$ws.Range("C25").Select()

# --- Column C width widened slightly to fit the new notes ---
$ws.Columns("C").ColumnWidth = 25.166666666666668
